$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mark bug #17 ("Running continues after game ended") as Fixed by Sandro
$ws.Range("B17").Value = "Fixed"
$ws.Range("C17").Value = "Sandro"

# Copy the "fixed" diagonal-strike formatting from another fixed row onto A17
$ws.Range("A2").Copy()
$ws.Range("A17").PasteSpecial(-4122) # xlPasteFormats

# Update the active selection to A12
$ws.Range("A12").Select()
